$d = $word.ActiveDocument
$t = $d.Tables(1)

# Map of 1-based table row -> new cell text.
# (Row indices correspond to the single-column table rows in document order.)
$changes = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "286"
    5  = "0.00001"
    6  = "0.00063"
    7  = "0.00011"
    9  = "0.00016"
    10 = "0.00018"
    11 = "0.00022"
    12 = "0.03480"
    44 = "99.88"
    45 = "0.03"
    46 = "29"
}

foreach ($rowIndex in $changes.Keys) {
    $cell = $t.Cell($rowIndex, 1)
    $cell.Range.Text = $changes[$rowIndex]
}
